$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Sprint 1 Evaluation - first bullet: extend sentence with Mitchell/James info
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "Both tasks set out for the first sprint were completed.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Both tasks set out for the first sprint were completed. Mitchell completed the basic front end and James investigated potential methods for the 2FA.",
    2)

# ---------------------------------------------------------------------------
# 2. Sprint 2 section reorganisation
# ---------------------------------------------------------------------------

# 2a. Insert a new bullet ("It was also discussed...") right after the
#     "For the second sprint..." planning bullet.
$r = $d.Content
$r.Find.Execute(
    "For the second sprint we decided to get the two-factor authenticator working separately from the main application to discover how it works and how it would be implemented in the bank application.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.InsertParagraphAfter()
    $newPara = $r.Paragraphs(1).Next()
    $newRange = $newPara.Range
    $newRange.Text = "It was also discussed that further detail should be added to the User Stories and the SRS document, such as the acceptance criteria of each story and the requirements for the project."

    # 2b. Insert a new "Sprint Evaluation:" Heading2 paragraph right after the
    #     paragraph we just created (the heading is being relocated here).
    $newRange.InsertParagraphAfter()
    $headingPara = $newRange.Paragraphs(1).Next()
    $headingRange = $headingPara.Range
    $headingRange.Text = "Sprint Evaluation:"
    $headingRange.Style = "Heading 2"
}

# 2c. Modify the "two-factor authenticator ... purpose built project." bullet.
$r = $d.Content
$r.Find.Execute(
    "The two-factor authenticator was successfully implemented on a separate purpose built project.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The two-factor authenticator was successfully implemented by James a separate purpose-built project.",
    2)

# 2d. Delete the old "Sprint Evaluation:" heading that used to sit directly
#     before the bullet we just modified (it has been relocated in step 2b).
$r = $d.Content
$r.Find.Execute(
    "The two-factor authenticator was successfully implemented by James a separate purpose-built project.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $para = $r.Paragraphs(1)
    $prevPara = $para.Previous()
    if ($prevPara.Range.Text -match "Sprint Evaluation:") {
        $prevPara.Range.Delete()
    }
}

# 2e. Insert a new bullet ("Michael also updated...") right after the
#     "testing authenticator" bullet.
$r = $d.Content
$r.Find.Execute(
    "The testing authenticator was not pushed to the main branch as it was not implemented to the main application.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.InsertParagraphAfter()
    $newPara = $r.Paragraphs(1).Next()
    $newRange = $newPara.Range
    $newRange.Text = "Michael also updated the User Stories to contain the acceptance criteria and further revised the SRS to allow for a deeper understanding of the project."
}

# ---------------------------------------------------------------------------
# 3. Sprint 3 Evaluation - database encryption bullet
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "Using the previously researched methods on database encryption a database was set up and the data stored on it was encrypted.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Using the previously researched methods on database encryption a database was set up by Mitchell and the data stored on it was encrypted.",
    2)

# ---------------------------------------------------------------------------
# 4. Sprint 4 Evaluation - transaction history / accounts balance bullets
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "The transaction history was added to the main menu.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The transaction history was added to the main menu by Mitchell and Michael.",
    2)

$r = $d.Content
$r.Find.Execute(
    "The accounts balance was added to the main menu.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The accounts balance was added to the main menu by Mitchell and Michael.",
    2)

# ---------------------------------------------------------------------------
# 5. Sprint 5 Evaluation - final bullet
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "All features were finished and pushed onto the main branch of the GitHub.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "All features were finished and pushed onto the main branch of the GitHub by Mitchell and James.",
    2)
